$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toni Alimi's talk: fill in the real title + abstract (were placeholders: "Title
# Coming Soon!" / " ")
$ws.Range("G29").Value = "Justifying Slavery: The Intellectual Background of the Reconstruction Amendments"
$ws.Range("H29").Value = " Recent appraisals of the Thirteenth Amendment to the United States Constitution often note with alarm that slavery remains a legal form of punishment and lament for how this loophole has been weaponized against black Americans. I’m interested in how this alarm reflects an Aristotelian attitude towards slavery (more on what that means in the talk!), and about how an intellectual history of various justifications for slavery can help us understand what’s going on in the Thirteenth Amendment."

# Small copyedit to Donald Rakow's abstract: hyphen -> em dash
$ws.Range("H28").Value = "Do you find yourself freaking out over the state of the global environment, climate change, and loss of biodiversity? If so, you’re not alone — eco-anxiety has become one of the leading mental health issues in the third decade of the 21st century. Public gardens, which include botanic gardens, arboreta, conservatories, and historic landscapes, are addressing these seemingly overwhelming environmental challenges through their extensive collections, protection of natural areas, preservation of endangered species, and through programming and messaging. Will public gardens single-handedly solve all of these problems? No, but they can be a significant contributor to the solutions. In this talk you’ll learn some of the ways they are making a difference."

# Leave the cursor where the author left it when they saved
$ws.Activate()
$ws.Range("F33").Select()
